# Apply "automatic update of files" changes:
#  1. Column C (Förändrad) for every data row (2-18): 45184 -> 45186
#  2. For rows 2 and 3, the HYPERLINK() formulas in columns S, T, V, W, X, Y
#     gain a second argument: the "friendly name" (the value of column A,
#     e.g. "A 13326-2020"), so the cell shows that text instead of the URL.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the "Förändrad" date in column C for rows 2 through 18 ---
for ($row = 2; $row -le 18; $row++) {
    $ws.Cells.Item($row, 3).Value = 45186
}

# --- 2. Add friendly-name argument to HYPERLINK formulas on rows 2 and 3 ---
$hyperlinkCols = @("S", "T", "V", "W", "X", "Y")

$linkInfo = @{
    2 = @{
        Name = "A 13326-2020"
        S = 'HYPERLINK("https://klasma.github.io/Logging_HUDDINGE/artfynd/A 13326-2020.xlsx", "A 13326-2020")'
        T = 'HYPERLINK("https://klasma.github.io/Logging_HUDDINGE/kartor/A 13326-2020.png", "A 13326-2020")'
        V = 'HYPERLINK("https://klasma.github.io/Logging_HUDDINGE/klagomål/A 13326-2020.docx", "A 13326-2020")'
        W = 'HYPERLINK("https://klasma.github.io/Logging_HUDDINGE/klagomålsmail/A 13326-2020.docx", "A 13326-2020")'
        X = 'HYPERLINK("https://klasma.github.io/Logging_HUDDINGE/tillsyn/A 13326-2020.docx", "A 13326-2020")'
        Y = 'HYPERLINK("https://klasma.github.io/Logging_HUDDINGE/tillsynsmail/A 13326-2020.docx", "A 13326-2020")'
    }
    3 = @{
        Name = "A 24233-2023"
        S = 'HYPERLINK("https://klasma.github.io/Logging_HUDDINGE/artfynd/A 24233-2023.xlsx", "A 24233-2023")'
        T = 'HYPERLINK("https://klasma.github.io/Logging_HUDDINGE/kartor/A 24233-2023.png", "A 24233-2023")'
        V = 'HYPERLINK("https://klasma.github.io/Logging_HUDDINGE/klagomål/A 24233-2023.docx", "A 24233-2023")'
        W = 'HYPERLINK("https://klasma.github.io/Logging_HUDDINGE/klagomålsmail/A 24233-2023.docx", "A 24233-2023")'
        X = 'HYPERLINK("https://klasma.github.io/Logging_HUDDINGE/tillsyn/A 24233-2023.docx", "A 24233-2023")'
        Y = 'HYPERLINK("https://klasma.github.io/Logging_HUDDINGE/tillsynsmail/A 24233-2023.docx", "A 24233-2023")'
    }
}

foreach ($row in $linkInfo.Keys) {
    $cols = $linkInfo[$row]
    foreach ($col in $hyperlinkCols) {
        $ws.Range("$col$row").Formula = "=" + $cols[$col]
    }
}
